$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "'28.103.02"
    "E2" = "'  -1.70%  "
    "D3" = "'1.895.25"
    "E3" = "'  -0.96%  "
    "E4" = "'  +0.08%  "
    "D5" = "'314.15"
    "E5" = "'  -0.25%  "
    "E6" = "'  +0.10%  "
    "D7" = "'0.5033"
    "E7" = "'  -0.42%  "
    "E8" = "'  -1.26%  "
    "D9" = "'0.09250"
    "E10" = "'  -2.67%  "
    "D11" = "'41.81"
    "E11" = "'  -1.12%  "
    "E12" = "'  -2.51%  "
    "E13" = "'  -1.88%  "
    "D14" = "'1.895.35"
    "E14" = "'  -1.31%  "
    "D15" = "'7.292"
    "E15" = "'  -3.79%  "
    "E16" = "'  +0.05%  "
    "D17" = "'92.48"
    "E17" = "'  -1.75%  "
    "D18" = "'0.00001109"
    "E18" = "'  -2.83%  "
    "D19" = "'0.06655"
    "E19" = "'  -0.05%  "
    "E20" = "'  -1.58%  "
    "E21" = "'  +0.05%  "
    "D22" = "'6.209"
    "E22" = "'  -1.64%  "
    "D23" = "'28.157.73"
    "E23" = "'  -1.71%  "
    "E24" = "'  -0.25%  "
    "E25" = "'  +1.84%  "
    "D26" = "'2.124.53"
    "E26" = "'  -0.53%  "
    "D27" = "'2.547"
    "E27" = "'  -7.25%  "
    "E28" = "'  -2.03%  "
    "D29" = "'158.26"
    "E29" = "'  -0.65%  "
    "D30" = "'127.02"
    "E30" = "'  -1.48%  "
    "E31" = "'  -2.38%  "
    "E32" = "'  -1.55%  "
    "D33" = "'5.617"
    "E33" = "'  -2.38%  "
    "E34" = "'  -0.94%  "
    "D35" = "'9.593"
    "E35" = "'  -2.70%  "
    "D36" = "'1.366"
    "E36" = "'  +14.98%  "
    "D37" = "'0.06608"
    "E37" = "'  -2.92%  "
    "D38" = "'0.02405"
    "E38" = "'  -1.54%  "
    "D39" = "'0.2209"
    "E39" = "'  -1.06%  "
    "D40" = "'1.224"
    "E40" = "'  -3.94%  "
    "D41" = "'0.6484"
    "E41" = "'  +0.61%  "
    "D42" = "'11.44"
    "E42" = "'  -2.97%  "
    "D43" = "'4.974"
    "E43" = "'  -2.47%  "
    "E44" = "'  +0.01%  "
    "B45" = "'EnergySwap"
    "C45" = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D45" = "'13.40"
    "E45" = "'  -2.37%  "
    "B46" = "'Decentraland"
    "C46" = "'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D46" = "'0.6087"
    "E46" = "'  -0.07%  "
    "D47" = "'1.304"
    "E47" = "'  +1.55%  "
    "D48" = "'3.688"
    "E48" = "'  -2.65%  "
    "E49" = "'  -2.12%  "
    "D50" = "'122.20"
    "E50" = "'  -2.38%  "
    "E51" = "'  -1.36%  "
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
